# Update the frozen-pane selection on the reduced_model_set_results sheet
$wb = $excel.ActiveWorkbook

$wsReduced = $wb.Worksheets.Item("reduced_model_set_results")
$wsReduced.Activate() | Out-Null
$wsReduced.Range("T2").Select() | Out-Null

# Switch to the results tracker sheet and fill in region 8 (starter_08) results
$ws = $wb.Worksheets.Item("model_rerun_results_tracker")
$ws.Activate() | Out-Null

$ws.Range("H450").Value = "NA"
$ws.Range("I450").Value = "NA"
$ws.Range("J450").Value = "NA"
$ws.Range("K450").Value = "T"
$ws.Range("L450").Value = "NA"
$ws.Range("H451").Value = "NA"
$ws.Range("I451").Value = "NA"
$ws.Range("J451").Value = "NA"
$ws.Range("K451").Value = "T"
$ws.Range("L451").Value = "NA"
$ws.Range("H452").Value = "NA"
$ws.Range("I452").Value = "NA"
$ws.Range("J452").Value = "NA"
$ws.Range("K452").Value = "T"
$ws.Range("L452").Value = "NA"
$ws.Range("H453").Value = "NA"
$ws.Range("I453").Value = "NA"
$ws.Range("J453").Value = "NA"
$ws.Range("K453").Value = "T"
$ws.Range("L453").Value = "NA"
$ws.Range("H454").Value = "NA"
$ws.Range("I454").Value = "NA"
$ws.Range("J454").Value = "NA"
$ws.Range("K454").Value = "T"
$ws.Range("L454").Value = "NA"
$ws.Range("H455").Value = "NA"
$ws.Range("I455").Value = "NA"
$ws.Range("J455").Value = "NA"
$ws.Range("K455").Value = "T"
$ws.Range("L455").Value = "NA"
$ws.Range("H456").Value = "NA"
$ws.Range("I456").Value = "NA"
$ws.Range("J456").Value = "NA"
$ws.Range("K456").Value = "T"
$ws.Range("L456").Value = "NA"
$ws.Range("H457").Value = "NA"
$ws.Range("I457").Value = "NA"
$ws.Range("J457").Value = "NA"
$ws.Range("K457").Value = "T"
$ws.Range("L457").Value = "NA"
$ws.Range("H458").Value = 12
$ws.Range("I458").Value = 3073
$ws.Range("J458").Value = 1599.48
$ws.Range("K458").Value = "T"
$ws.Range("L458").Value = "T"
$ws.Range("H459").Value = 13
$ws.Range("I459").Value = 2691
$ws.Range("J459").Value = 1275.74
$ws.Range("K459").Value = "T"
$ws.Range("L459").Value = "T"
$ws.Range("H460").Value = 12
$ws.Range("I460").Value = 3020
$ws.Range("J460").Value = 1495.68
$ws.Range("K460").Value = "T"
$ws.Range("L460").Value = "T"
$ws.Range("H461").Value = 12
$ws.Range("I461").Value = 2965
$ws.Range("J461").Value = 1454.82
$ws.Range("K461").Value = "T"
$ws.Range("L461").Value = "T"
$ws.Range("H462").Value = 12
$ws.Range("I462").Value = 3073
$ws.Range("J462").Value = 1599.48
$ws.Range("K462").Value = "T"
$ws.Range("L462").Value = "T"
$ws.Range("H463").Value = 13
$ws.Range("I463").Value = 2691
$ws.Range("J463").Value = 1275.71
$ws.Range("K463").Value = "T"
$ws.Range("L463").Value = "T"
$ws.Range("H464").Value = 12
$ws.Range("I464").Value = 3116
$ws.Range("J464").Value = 1585.86
$ws.Range("K464").Value = "T"
$ws.Range("L464").Value = "T"
$ws.Range("H465").Value = 12
$ws.Range("I465").Value = 2965
$ws.Range("J465").Value = 1454.82
$ws.Range("K465").Value = "T"
$ws.Range("L465").Value = "T"
$ws.Range("H466").Value = 5
$ws.Range("I466").Value = 2277
$ws.Range("J466").Value = 821.07
$ws.Range("K466").Value = "T"
$ws.Range("L466").Value = "T"
$ws.Range("H467").Value = 5
$ws.Range("I467").Value = 2245
$ws.Range("J467").Value = 790.34
$ws.Range("K467").Value = "T"
$ws.Range("L467").Value = "T"
$ws.Range("H468").Value = 6
$ws.Range("I468").Value = 2394
$ws.Range("J468").Value = 889.45
$ws.Range("K468").Value = "T"
$ws.Range("L468").Value = "T"
$ws.Range("H469").Value = 6
$ws.Range("I469").Value = 2456
$ws.Range("J469").Value = 985.54
$ws.Range("K469").Value = "T"
$ws.Range("L469").Value = "T"
$ws.Range("H470").Value = 4
$ws.Range("I470").Value = 2190
$ws.Range("J470").Value = 774.24
$ws.Range("K470").Value = "T"
$ws.Range("L470").Value = "T"
$ws.Range("H471").Value = 5
$ws.Range("I471").Value = 2447
$ws.Range("J471").Value = 949.89
$ws.Range("K471").Value = "T"
$ws.Range("L471").Value = "T"
$ws.Range("H472").Value = 5
$ws.Range("I472").Value = 2258
$ws.Range("J472").Value = 840.19
$ws.Range("K472").Value = "T"
$ws.Range("L472").Value = "T"
$ws.Range("H473").Value = 4
$ws.Range("I473").Value = 2132
$ws.Range("J473").Value = 721.2
$ws.Range("K473").Value = "T"
$ws.Range("L473").Value = "T"
$ws.Range("H474").Value = 5
$ws.Range("I474").Value = 2296
$ws.Range("J474").Value = 832.52
$ws.Range("K474").Value = "T"
$ws.Range("L474").Value = "T"
$ws.Range("H475").Value = 5
$ws.Range("I475").Value = 2228
$ws.Range("J475").Value = 783.03
$ws.Range("K475").Value = "T"
$ws.Range("L475").Value = "T"
$ws.Range("H476").Value = 6
$ws.Range("I476").Value = 2259
$ws.Range("J476").Value = 854.93
$ws.Range("K476").Value = "T"
$ws.Range("L476").Value = "T"
$ws.Range("H477").Value = 6
$ws.Range("I477").Value = 2402
$ws.Range("J477").Value = 935.61
$ws.Range("K477").Value = "T"
$ws.Range("L477").Value = "T"
$ws.Range("H478").Value = 4
$ws.Range("I478").Value = 2175
$ws.Range("J478").Value = 759.52
$ws.Range("K478").Value = "T"
$ws.Range("L478").Value = "T"
$ws.Range("H479").Value = 4
$ws.Range("I479").Value = 2249
$ws.Range("J479").Value = 822.85
$ws.Range("K479").Value = "T"
$ws.Range("L479").Value = "T"
$ws.Range("H480").Value = 4
$ws.Range("I480").Value = 2259
$ws.Range("J480").Value = 799.63
$ws.Range("K480").Value = "T"
$ws.Range("L480").Value = "T"
$ws.Range("H481").Value = 4
$ws.Range("I481").Value = 2189
$ws.Range("J481").Value = 761.56
$ws.Range("K481").Value = "T"
$ws.Range("L481").Value = "T"

# Restore the active cell / view position as left by the editor
$ws.Range("C1").Select() | Out-Null
